$wb = $excel.ActiveWorkbook

# --- Sheet "linear" ---
$wsLinear = $wb.Worksheets.Item("linear")
$wsLinear.Range("B2").Value = 0.004105376846097051
$wsLinear.Range("B3").Value = -0.07898817841423227
$wsLinear.Range("B4").Value = 1.353521011482407

# --- Sheet "non-linear" ---
$wsNonLinear = $wb.Worksheets.Item("non-linear")
$wsNonLinear.Range("B2").Value = 0.02052115074596907
$wsNonLinear.Range("B3").Value = 0.01509938526654823
$wsNonLinear.Range("B4").Value = 1.379653801527025
$wsNonLinear.Range("B5").Value = 0.08066323695490685
$wsNonLinear.Range("B6").Value = -0.276795127221192
$wsNonLinear.Range("B7").Value = 1.324415593620945
